$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Fix the "CasesTab" Neo4j query (B2): the original query incorrectly
# returned an extra `Cohort` column (via an OPTIONAL MATCH on (co:cohort))
# that does not belong in the Cases export. Trim the trailing
# `coalesce(co.cohort_description, '') AS `Cohort`` line (and the now
# dangling trailing comma/newline on the previous line) so the query ends
# cleanly on the `Response to Treatment` column.
# ---------------------------------------------------------------------------
$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Irish Wolfhound'] 
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesQuery

# ---------------------------------------------------------------------------
# Excel re-wraps/re-flows the three long query cells after the edit; match
# the resulting row heights (content shrank from 19 to 17 wrapped lines for
# row 2, other rows keep their existing wrapped-line counts).
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# ---------------------------------------------------------------------------
# Move the active selection / scrolled view back up to the (now shorter)
# CasesTab row instead of the FilesTab row.
# ---------------------------------------------------------------------------
$ws.Range("B2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
